$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.058, -0.155, -0.156, 152, -46, -110),
    @(-0.09,  0.101,  0.323,  43,  -45,  118),
    @(-0.52,  0.101,  0.018,  158, -38, -150),
    @(-0.134, 0.13,   0.051, -108,  56,   97),
    @(0.329,  0.214, -0.233, -42,   25,   83),
    @(0.059, -0.108,  0.133, -107,  34, -105),
    @(0.113,  0.162,  0.193,  106, -42,  173),
    @(0.304, -0.126,  0.412,  168, -41,  -30),
    @(0.168, -0.039, -0.207,  -60,  36,  -78),
    @(0.234,  0.027,  0.532,   89, -13,   -1)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}
